$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Quarter 1 "Print" expense from 10000 to 20000.
# This triggers automatic recalculation of the dependent formulas/values
# on the sheet (D7:D14, D19, D27, D33, D41, D48, C49, etc.)
$ws.Range("C8").Value = 20000

# Update the active cell selection to K8 (as reflected in the saved view state).
$ws.Range("K8").Select()
